$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1656050955414013
$ws.Range("C2").Value = 0.6178343949044586
$ws.Range("J2").Value = 0.009554140127388535
$ws.Range("P2").Value = 0.1273885350318471
$ws.Range("S2").Value = 0.07961783439490445
$ws.Range("C3").Value = 0.02941176470588235
$ws.Range("J3").Value = 0.02450980392156863
$ws.Range("P3").Value = 0.7843137254901961
$ws.Range("S3").Value = 0.1617647058823529
$ws.Range("J4").Value = 0.03773584905660377
$ws.Range("P4").Value = 0.7735849056603774
$ws.Range("S4").Value = 0.1886792452830189
$ws.Range("B6").Value = 0.0603448275862069
$ws.Range("D6").Value = 0.01293103448275862
$ws.Range("F6").Value = 0.08189655172413793
$ws.Range("J6").Value = 0.2543103448275862
$ws.Range("O6").Value = 0.004310344827586207
$ws.Range("Q6").Value = 0.1163793103448276
$ws.Range("R6").Value = 0.08189655172413793
$ws.Range("S6").Value = 0.3879310344827586
$ws.Range("B7").Value = 0.08620689655172414
$ws.Range("D7").Value = 0.02873563218390805
$ws.Range("F7").Value = 0.06896551724137931
$ws.Range("J7").Value = 0.09195402298850575
$ws.Range("O7").Value = 0.005747126436781609
$ws.Range("Q7").Value = 0.1724137931034483
$ws.Range("R7").Value = 0.07471264367816093
$ws.Range("S7").Value = 0.4712643678160919
$ws.Range("B8").Value = 0.09394572025052192
$ws.Range("D8").Value = 0.0208768267223382
$ws.Range("F8").Value = 0.08559498956158663
$ws.Range("J8").Value = 0.1106471816283925
$ws.Range("O8").Value = 0.0208768267223382
$ws.Range("Q8").Value = 0.1419624217118998
$ws.Range("R8").Value = 0.1064718162839248
$ws.Range("S8").Value = 0.4196242171189979
$ws.Range("B9").Value = 0.1473684210526316
$ws.Range("D9").Value = 0.02105263157894737
$ws.Range("F9").Value = 0.07894736842105263
$ws.Range("J9").Value = 0.1052631578947368
$ws.Range("O9").Value = 0.01052631578947368
$ws.Range("Q9").Value = 0.1842105263157895
$ws.Range("R9").Value = 0.08947368421052632
$ws.Range("S9").Value = 0.3631578947368421
$ws.Range("B10").Value = 0.1267942583732057
$ws.Range("D10").Value = 0.02472089314194577
$ws.Range("E10").Value = 0.0007974481658692185
$ws.Range("F10").Value = 0.06618819776714513
$ws.Range("J10").Value = 0.1108452950558214
$ws.Range("O10").Value = 0.0127591706539075
$ws.Range("Q10").Value = 0.196969696969697
$ws.Range("R10").Value = 0.09250398724082935
$ws.Range("S10").Value = 0.3684210526315789
$ws.Range("G11").Value = 0.1433962264150943
$ws.Range("J11").Value = 0.1056603773584906
$ws.Range("K11").Value = 0.2037735849056604
$ws.Range("L11").Value = 0.5283018867924528
$ws.Range("S11").Value = 0.01886792452830189
$ws.Range("G12").Value = 0.7638888888888888
$ws.Range("J12").Value = 0.1597222222222222
$ws.Range("K12").Value = 0.006944444444444444
$ws.Range("L12").Value = 0.03472222222222222
$ws.Range("S12").Value = 0.03472222222222222
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2916666666666667
$ws.Range("S13").Value = 0.04166666666666666
$ws.Range("F15").Value = 0.01507537688442211
$ws.Range("H15").Value = 0.2311557788944724
$ws.Range("I15").Value = 0.07537688442211055
$ws.Range("J15").Value = 0.321608040201005
$ws.Range("K15").Value = 0.04522613065326633
$ws.Range("M15").Value = 0.01507537688442211
$ws.Range("O15").Value = 0.04522613065326633
$ws.Range("S15").Value = 0.2512562814070352
$ws.Range("F16").Value = 0.02575107296137339
$ws.Range("H16").Value = 0.1716738197424893
$ws.Range("I16").Value = 0.07725321888412018
$ws.Range("J16").Value = 0.4120171673819742
$ws.Range("K16").Value = 0.1072961373390558
$ws.Range("M16").Value = 0.03004291845493562
$ws.Range("O16").Value = 0.04721030042918455
$ws.Range("S16").Value = 0.1287553648068669
$ws.Range("F17").Value = 0.0196078431372549
$ws.Range("H17").Value = 0.1593137254901961
$ws.Range("I17").Value = 0.09313725490196079
$ws.Range("J17").Value = 0.4558823529411765
$ws.Range("K17").Value = 0.07107843137254902
$ws.Range("M17").Value = 0.0196078431372549
$ws.Range("N17").Value = 0.002450980392156863
$ws.Range("O17").Value = 0.07598039215686274
$ws.Range("S17").Value = 0.1029411764705882
$ws.Range("F18").Value = 0.02347417840375587
$ws.Range("H18").Value = 0.1830985915492958
$ws.Range("I18").Value = 0.07981220657276995
$ws.Range("J18").Value = 0.4178403755868544
$ws.Range("K18").Value = 0.0892018779342723
$ws.Range("M18").Value = 0.009389671361502348
$ws.Range("O18").Value = 0.08450704225352113
$ws.Range("S18").Value = 0.1126760563380282
$ws.Range("F19").Value = 0.01352426412092283
$ws.Range("H19").Value = 0.2330946698488465
$ws.Range("I19").Value = 0.08035003977724742
$ws.Range("J19").Value = 0.3794749403341289
$ws.Range("K19").Value = 0.09785202863961814
$ws.Range("M19").Value = 0.02307080350039777
$ws.Range("N19").Value = 0.0007955449482895784
$ws.Range("O19").Value = 0.05807478122513922
$ws.Range("S19").Value = 0.1137629276054097
